$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New roster data (player, position, team) replacing the old A2:C19 block.
$data = @(
    @("Jalen Brunson", "PG", "New York Knicks"),
    @("Immanuel Quickley", "PG,SG", "Toronto Raptors"),
    @("Max Christie", "SG,SF", "Dallas Mavericks"),
    @("Desmond Bane", "SG,SF", "Memphis Grizzlies"),
    @("Jimmy Butler III", "SF,PF", "Golden State Warriors"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Jeremy Sochan", "SF,PF", "San Antonio Spurs"),
    @("Jabari Smith Jr.", "PF,C", "Houston Rockets"),
    @("Walker Kessler", "C", "Utah Jazz"),
    @("Myles Turner", "C", "Indiana Pacers"),
    @("Kawhi Leonard", "SG,SF,PF", "LA Clippers"),
    @("LeBron James", "SF,PF", "Los Angeles Lakers"),
    @("Christian Braun", "SG,SF", "Denver Nuggets"),
    @("Trae Young", "PG", "Atlanta Hawks"),
    @("Devin Booker", "PG,SG", "Phoenix Suns"),
    @("Brandon Ingram", "SG,SF,PF", "Toronto Raptors"),
    @("Norman Powell", "SG,SF", "LA Clippers")
)

# Old table went through row 19; new table only goes through row 18, so
# clear the old extent first, then write the new values.
$ws.Range("A2:C19").ClearContents()

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
